# Add a new paragraph "V2" (red, FF0000) right after the existing "V1" paragraph.

$d = $word.ActiveDocument

# Move to the very end of the document content and append a new paragraph.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# Move into the freshly created (empty) paragraph and insert its text.
$end.Collapse(0)
$end.InsertAfter("V2")

# Color the new paragraph's text (and paragraph mark) red (RGB FF0000 -> 255).
$d.Paragraphs($d.Paragraphs.Count).Range.Font.Color = 255
